# Append match rows 80..86 (sheet rows 81..87) to the Ecuador Liga Pro 2023
# results sheet, mirroring the source scraper's "Atualizado por script" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: sheetRow, Indice, pais, torneio, temporada, data_partida,
# home, home_ft_gols, away, away_ft_gols,
# home_opening_odds, home_opening_data_hora, home_closing_odds, home_closing_data_hora,
# draw_opening_odds, draw_opening_data_hora, draw_closing_odds, draw_closing_data_hora,
# away_opening_odds, away_opening_data_hora, away_closing_odds, away_closing_data_hora,
# url_partida
$data = @(
    @(81, 80, "ecuador", "liga-pro", "2023", 45227.08333333334, "Libertad", 3, "Gualaceo", 0, 1.74, "22/10/2023 20:15", 2.21, "28/10/2023 01:52", 3.66, "22/10/2023 20:15", 3.15, "28/10/2023 01:54", 4.76, "22/10/2023 20:15", 3.64, "28/10/2023 01:54", "https://www.betexplorer.com/football/ecuador/liga-pro/libertad-gualaceo/CMSKLlLb/"),
    @(82, 81, "ecuador", "liga-pro", "2023", 45227.83333333334, "Ind. del Valle", 2, "Dep. Cuenca", 0, 1.36, "22/10/2023 22:13", 1.36, "28/10/2023 08:37", 4.74, "22/10/2023 22:13", 4.93, "28/10/2023 19:38", 7.22, "22/10/2023 22:13", 9.050000000000001, "28/10/2023 19:38", "https://www.betexplorer.com/football/ecuador/liga-pro/independiente-del-valle-dep-cuenca/YVTGMU5h/"),
    @(83, 82, "ecuador", "liga-pro", "2023", 45228.04166666666, "Emelec", 2, "U. Catolica", 1, 2.12, "24/10/2023 02:12", 2.23, "29/10/2023 00:56", 3.51, "24/10/2023 02:12", 3.48, "29/10/2023 00:56", 3.18, "24/10/2023 02:12", 3.23, "29/10/2023 00:56", "https://www.betexplorer.com/football/ecuador/liga-pro/emelec-u-catolica/Mi42375U/"),
    @(84, 83, "ecuador", "liga-pro", "2023", 45228.79166666666, "Delfin", 4, "Cumbaya", 0, 1.65, "24/10/2023 02:12", 1.49, "29/10/2023 18:51", 3.66, "24/10/2023 02:12", 4.01, "29/10/2023 18:51", 5.6, "24/10/2023 02:12", 7.64, "29/10/2023 18:51", "https://www.betexplorer.com/football/ecuador/liga-pro/delfin-cumbaya/nFV8Ojyt/"),
    @(85, 84, "ecuador", "liga-pro", "2023", 45228.89583333334, "Aucas", 4, "Mushuc Runa", 0, 1.69, "22/10/2023 22:42", 1.8, "29/10/2023 21:26", 3.74, "22/10/2023 22:42", 3.68, "29/10/2023 21:26", 4.59, "22/10/2023 22:42", 4.54, "29/10/2023 21:26", "https://www.betexplorer.com/football/ecuador/liga-pro/aucas-mushuc-runa/Q5UCNAjn/"),
    @(86, 85, "ecuador", "liga-pro", "2023", 45229, "Orense", 2, "Barcelona SC", 3, 2.8, "23/10/2023 01:12", 2.62, "29/10/2023 23:55", 3.13, "23/10/2023 01:12", 3.23, "29/10/2023 23:25", 2.54, "23/10/2023 01:12", 2.83, "29/10/2023 23:55", "https://www.betexplorer.com/football/ecuador/liga-pro/orense-barcelona-sc/zy2j69LB/"),
    @(87, 86, "ecuador", "liga-pro", "2023", 45230.04166666666, "EL Nacional", 2, "Guayaquil City", 0, 1.42, "24/10/2023 02:12", 1.41, "31/10/2023 00:56", 4.71, "24/10/2023 02:12", 4.63, "31/10/2023 00:56", 6.93, "24/10/2023 02:12", 5.94, "31/10/2023 00:56", "https://www.betexplorer.com/football/ecuador/liga-pro/el-nacional-guayaquil-city/ER0f5TzI/")
)

$firstNewRow = 81
$lastNewRow = 87

# Replicate the existing column formatting (bold/bordered index column A,
# datetime-formatted column E) onto the freshly appended rows before writing
# values, exactly like the existing rows 2..80.
$ws.Range("A2").Copy()
$ws.Range("A" + $firstNewRow + ":A" + $lastNewRow).PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E" + $firstNewRow + ":E" + $lastNewRow).PasteSpecial(-4122)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    # Column D ("temporada") is stored as literal text ("2023"), not a
    # number. A direct .Value assignment of a digit-only string gets
    # auto-coerced to a numeric cell (same "smart" parsing Excel applies
    # when you type a bare number into a cell), so stage it through a
    # text formula and flatten it back to a literal value below instead.
    $ws.Cells.Item($r, 4).Formula = '="' + $row[4] + '"'
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
    $ws.Cells.Item($r, 19).Value = $row[19]
    $ws.Cells.Item($r, 20).Value = $row[20]
    $ws.Cells.Item($r, 21).Value = $row[21]
    $ws.Cells.Item($r, 22).Value = $row[22]
}

# Flatten the column-D helper formulas down to plain text literals (no
# leftover <f> element), matching the inline-string cells used by every
# other row in this column.
$dRange = $ws.Range("D" + $firstNewRow + ":D" + $lastNewRow)
$dRange.Copy()
$dRange.PasteSpecial(-4163)
